$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (N_Calib_1=20, N_Calib_2=40)
$ws.Range("C2").Value = -0.1526314041118104
$ws.Range("D2").Value = 0.8800801668149112

# Row 3 (N_Calib_1=20, N_Calib_2=60)
$ws.Range("C3").Value = 0.1907607838301029
$ws.Range("D3").Value = 0.8504624286690567

# Row 4 (N_Calib_1=20, N_Calib_2=100)
$ws.Range("C4").Value = 2.142644983921303
$ws.Range("D4").Value = 0.04345948721276871

# Row 5 (N_Calib_1=20, N_Calib_2=200)
$ws.Range("C5").Value = 2.218976492267076
$ws.Range("D5").Value = 0.03711273325222852

# Row 6 (N_Calib_1=40, N_Calib_2=60)
$ws.Range("C6").Value = 0.2975835928242627
$ws.Range("D6").Value = 0.7688129552271723

# Row 7 (N_Calib_1=40, N_Calib_2=100)
$ws.Range("C7").Value = 2.786975984990386
$ws.Range("D7").Value = 0.0107485615817795

# Row 8 (N_Calib_1=40, N_Calib_2=200)
$ws.Range("C8").Value = 2.615667005362166
$ws.Range("D8").Value = 0.01578742879189399
$ws.Range("G8").Value = "Sí"

# Row 9 (N_Calib_1=60, N_Calib_2=100)
$ws.Range("C9").Value = 1.86432088353592
$ws.Range("D9").Value = 0.07567538135986029
$ws.Range("G9").Value = "No"

# Row 10 (N_Calib_1=60, N_Calib_2=200)
$ws.Range("C10").Value = 2.487984036173935
$ws.Range("D10").Value = 0.02090818848079112
$ws.Range("G10").Value = "Sí"

# Row 11 (N_Calib_1=100, N_Calib_2=200)
$ws.Range("C11").Value = -0.4511468237515471
$ws.Range("D11").Value = 0.6562992614982801
